$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

$ws.Range("A34:K34").Value = "nan"
$ws.Range("M34").Value = "قطع سير 1270 نتيجه خلل ف عيار "
